# Applies the odds updates described in the diff to Sheet1, rows 3 and 4.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 updates
$ws.Range("G3").Value = 2.1
$ws.Range("I3").Value = 3.8

# Row 4 updates
$ws.Range("G4").Value = 2
$ws.Range("I4").Value = 3.6
$ws.Range("J4").Value = 2.6
$ws.Range("L4").Value = 4
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 12
$ws.Range("U4").Value = 1.75
$ws.Range("V4").Value = 2
$ws.Range("X4").Value = 9.5
$ws.Range("Y4").Value = 8.5
$ws.Range("Z4").Value = 17
$ws.Range("AA4").Value = 15
$ws.Range("AC4").Value = 12
$ws.Range("AD4").Value = 7
$ws.Range("AE4").Value = 15
$ws.Range("AF4").Value = 51
$ws.Range("AI4").Value = 19
$ws.Range("AL4").Value = 29
$ws.Range("AQ4").Value = 34
$ws.Range("AS4").Value = 126
$ws.Range("AX4").Value = 21
$ws.Range("AZ4").Value = 67
